$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Site glossary")

# Append two new glossary entries (Resource Abbreviation) after the
# existing last data row (row 103).
$ws.Range("A104").Value = "Resource Abbreviation"
$ws.Range("B104").Value = "PPTP"
$ws.Range("C104").Value = "Pediatric Preclinical Testing Program"

$ws.Range("A105").Value = "Resource Abbreviation"
$ws.Range("B105").Value = "DepMap"
$ws.Range("C105").Value = "Dependency Map"
